$d = $word.ActiveDocument

# Position a zero-length (collapsed) range at the very end of the
# document's main story, right before the closing sectPr, so the new
# OOXML is appended as sibling paragraphs rather than merged into the
# last existing paragraph.
$content = $d.Content
$endRange = $d.Range($content.End, $content.End)

# Four new paragraphs: a page-break-only paragraph, a blank paragraph,
# an "Asdasd" paragraph (with the spell-check proofErr markers Word
# leaves around a flagged word), and a "test" paragraph.
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:br w:type="page"/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Asdasd</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>test</w:t></w:r></w:p>
'@

$endRange.InsertXML($xml)
